# Update cryptos list cell values per source diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'70.147.49"
$ws.Range("E2").Value = "'  -1.18%  "

$ws.Range("D3").Value = "'3.580.46"
$ws.Range("E3").Value = "'  -2.00%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "'  -0.23%  "

$ws.Range("D5").Value = "'578.23"
$ws.Range("E5").Value = "'  -2.97%  "

$ws.Range("D6").Value = "'186.66"
$ws.Range("E6").Value = "'  -4.42%  "

$ws.Range("D7").Value = "'3.574.10"
$ws.Range("E7").Value = "'  -2.02%  "

$ws.Range("D8").Value = "'0.620"
$ws.Range("E8").Value = "'  -4.22%  "

$ws.Range("E9").Value = "'  -0.01%  "

$ws.Range("D10").Value = "'0.185"
$ws.Range("E10").Value = "'  +0.39%  "

$ws.Range("D11").Value = "'0.651"
$ws.Range("E11").Value = "'  -3.68%  "

$ws.Range("D12").Value = "'55.19"
$ws.Range("E12").Value = "'  -4.71%  "

$ws.Range("D13").Value = "'0.0000308"
$ws.Range("E13").Value = "'  +4.30%  "

$ws.Range("D14").Value = "'9.55"
$ws.Range("E14").Value = "'  -4.21%  "

$ws.Range("D15").Value = "'4.142.54"
$ws.Range("E15").Value = "'  -2.31%  "

$ws.Range("D16").Value = "'19.70"
$ws.Range("E16").Value = "'  -2.68%  "

$ws.Range("D17").Value = "'3.564.87"

$ws.Range("D18").Value = "'70.006.09"
$ws.Range("E18").Value = "'  -1.49%  "

$ws.Range("D19").Value = "'12.63"
$ws.Range("E19").Value = "'  -1.10%  "

$ws.Range("E20").Value = "'  -0.86%  "

$ws.Range("E21").Value = "'  -3.01%  "

$ws.Range("D22").Value = "'490.37"
$ws.Range("E22").Value = "'  +0.40%  "

$ws.Range("D23").Value = "'19.13"
$ws.Range("E23").Value = "'  +0.18%  "

$ws.Range("D24").Value = "'4.93"
$ws.Range("E24").Value = "'  -6.36%  "

$ws.Range("D25").Value = "'4.39"
$ws.Range("E25").Value = "'  -2.10%  "

$ws.Range("D26").Value = "'95.45"
$ws.Range("E26").Value = "'  +4.27%  "

$ws.Range("D27").Value = "'11.91"
$ws.Range("E27").Value = "'  +4.02%  "

$ws.Range("D28").Value = "'2.97"
$ws.Range("E28").Value = "'  -6.12%  "

$ws.Range("D29").Value = "'9.35"
$ws.Range("E29").Value = "'  -2.71%  "

$ws.Range("D30").Value = "'7.80"
$ws.Range("E30").Value = "'  -0.73%  "

$ws.Range("D31").Value = "'31.66"
$ws.Range("E31").Value = "'  -3.58%  "

$ws.Range("D32").Value = "'66.88"
$ws.Range("E32").Value = "'  +0.92%  "

$ws.Range("D33").Value = "'12.10"
$ws.Range("E33").Value = "'  -1.33%  "

$ws.Range("D34").Value = "'0.115"
$ws.Range("E34").Value = "'  -6.34%  "

$ws.Range("D35").Value = "'574.09"
$ws.Range("E35").Value = "'  -7.62%  "

$ws.Range("D36").Value = "'3.28"
$ws.Range("E36").Value = "'  +15.31%  "

$ws.Range("B37").Value = "'TheGraph"
$ws.Range("C37").Value = "'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D37").Value = "'0.413"
$ws.Range("E37").Value = "'  +0.26%  "

$ws.Range("B38").Value = "'InjectiveProtocol"
$ws.Range("C38").Value = "'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D38").Value = "'38.85"
$ws.Range("E38").Value = "'  -3.26%  "

$ws.Range("E39").Value = "'  +0.01%  "

$ws.Range("D40").Value = "'0.0₃0794"
$ws.Range("E40").Value = "'  -4.83%  "

$ws.Range("D41").Value = "'3.47"
$ws.Range("E41").Value = "'  -3.36%  "

$ws.Range("D42").Value = "'3.18"
$ws.Range("E42").Value = "'  -0.73%  "

$ws.Range("E43").Value = "'  -9.05%  "

$ws.Range("D44").Value = "'3.07"
$ws.Range("E44").Value = "'  -2.79%  "

$ws.Range("D45").Value = "'3.231.96"
$ws.Range("E45").Value = "'  -3.04%  "

$ws.Range("D46").Value = "'0.0444"
$ws.Range("E46").Value = "'  -2.57%  "

$ws.Range("D47").Value = "'3.45"
$ws.Range("E47").Value = "'  +3.96%  "

$ws.Range("D48").Value = "'9.61"
$ws.Range("E48").Value = "'  -0.33%  "

$ws.Range("B49").Value = "'OceanProtocol"
$ws.Range("C49").Value = "'https://coinranking.com/coin/aAKLSV5-0+oceanprotocol-ocean"
$ws.Range("D49").Value = "'1.58"
$ws.Range("E49").Value = "'  +31.20%  "

$ws.Range("B50").Value = "'Stellar"
$ws.Range("C50").Value = "'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D50").Value = "'0.136"
$ws.Range("E50").Value = "'  -2.26%  "

$ws.Range("D51").Value = "'0.997"
$ws.Range("E51").Value = "'  -0.34%  "
